# Apply the roster shuffle edit described by the diff.
# Rows 4, 5, 13, 14, 15, 16 get their Name / Position / Team values
# rewritten to reflect the new arrangement (a 6-row cycle of moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Jalen Suggs"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Orlando Magic"

$ws.Range("A5").Value = "Chris Paul"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "San Antonio Spurs"

$ws.Range("A13").Value = "Rudy Gobert"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Minnesota Timberwolves"

$ws.Range("A14").Value = "Klay Thompson"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Dallas Mavericks"

$ws.Range("A15").Value = "Jaylen Brown"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Boston Celtics"

$ws.Range("A16").Value = "Jakob Poeltl"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Toronto Raptors"
